$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 2.2
    "H2"  = 3.3
    "I2"  = 3.1
    "J2"  = 2.72
    "K2"  = 2.15
    "L2"  = 3.6
    "N2"  = 7.9
    "O2"  = 1.23
    "P2"  = 3.7
    "Q2"  = 1.72
    "R2"  = 2.05
    "S2"  = 1.36
    "T2"  = 2.92
    "U2"  = 1.55
    "V2"  = 2.3
    "W2"  = 9.75
    "Y2"  = 8.75
    "Z2"  = 23
    "AA2" = 16
    "AB2" = 21
    "AC2" = 7.9
    "AD2" = 6.5
    "AE2" = 11.5
    "AG2" = 250
    "AH2" = 10.75
    "AI2" = 17.5
    "AL2" = 25
    "AM2" = 28
    "AN2" = 4.35
    "AO2" = 11.25
    "AP2" = 16.5
    "AQ2" = 40
    "AR2" = 60
    "AT2" = 2.92
    "AU2" = 6.4
    "AX2" = 17
    "AY2" = 22
    "AZ2" = 80
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
